# Re-roll the randomized cue sequence ("lock in current version"):
# replace the word / image / category columns for every data row
# (row 1 is the header and is left untouched).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @(2,  "wenden",    "none",              "none"),
    @(3,  "schwimmen", "face/face017.jpg",  "face"),
    @(4,  "grüßen",    "house/house022.jpg","house"),
    @(5,  "rufen",     "none",              "none"),
    @(6,  "leiten",    "house/house028.jpg","house"),
    @(7,  "knien",     "house/house001.jpg","house"),
    @(8,  "enden",     "none",              "none"),
    @(9,  "decken",    "face/face001.jpg",  "face"),
    @(10, "töten",     "house/house023.jpg","house"),
    @(11, "opfern",    "none",              "none"),
    @(12, "runden",    "face/face021.jpg",  "face"),
    @(13, "planen",    "house/house010.jpg","house"),
    @(14, "schalten",  "none",              "none"),
    @(15, "triefen",   "face/face028.jpg",  "face"),
    @(16, "rühren",    "face/face008.jpg",  "face"),
    @(17, "klagen",    "none",              "none"),
    @(18, "helfen",    "face/face009.jpg",  "face"),
    @(19, "faulen",    "house/house008.jpg","house"),
    @(20, "sparen",    "none",              "none"),
    @(21, "kommen",    "house/house025.jpg","house"),
    @(22, "dringen",   "face/face014.jpg",  "face"),
    @(23, "hören",     "none",              "none"),
    @(24, "posten",    "face/face023.jpg",  "face"),
    @(25, "lehnen",    "face/face006.jpg",  "face"),
    @(26, "drohen",    "none",              "none"),
    @(27, "wohnen",    "house/house002.jpg","house"),
    @(28, "achten",    "face/face013.jpg",  "face"),
    @(29, "orten",     "none",              "none"),
    @(30, "sehen",     "house/house027.jpg","house"),
    @(31, "danken",    "house/house021.jpg","house"),
    @(32, "weigern",   "none",              "none"),
    @(33, "schütteln", "house/house019.jpg","house"),
    @(34, "segnen",    "house/house031.jpg","house"),
    @(35, "ächzen",    "none",              "none"),
    @(36, "regnen",    "house/house012.jpg","house"),
    @(37, "bilden",    "face/face027.jpg",  "face"),
    @(38, "kosten",    "none",              "none"),
    @(39, "öffnen",    "face/face000.jpg",  "face"),
    @(40, "fordern",   "face/face011.jpg",  "face"),
    @(41, "dauern",    "none",              "none"),
    @(42, "wundern",   "house/house029.jpg","house"),
    @(43, "erben",     "face/face025.jpg",  "face"),
    @(44, "stören",    "none",              "none"),
    @(45, "spielen",   "house/house026.jpg","house"),
    @(46, "brauchen",  "face/face003.jpg",  "face"),
    @(47, "bremsen",   "none",              "none"),
    @(48, "deuten",    "face/face019.jpg",  "face"),
    @(49, "kümmern",   "house/house005.jpg","house")
)

foreach ($row in $rows) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
}
